$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = "maa://21246 (91.41), maa://36684 (95.19), ***maa://22731 (6.67)"
$ws.Range("AF2").Value = "maa://25251 (92.79), ***maa://21730 (24.66), ***maa://39501 (20.83), *maa://36675 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (65.26), maa://20276 (86.13), *maa://22749 (72.73)"
$ws.Range("P3").Value = "maa://21249 (94.37), maa://26254 (96.43)"
$ws.Range("X3").Value = "maa://27396 (84.28), maa://27484 (96.49), maa://27480 (82.86)"
$ws.Range("D4").Value = "maa://24632 (93.79), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("X4").Value = "**maa://32495 (48.51), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (87.23)"
$ws.Range("AF4").Value = "*maa://30062 (63.27), ***maa://26209 (13.04), *maa://39394 (66.67)"
$ws.Range("D5").Value = "maa://21245 (84.21), maa://22744 (84.0)"
$ws.Range("A8").Value = "更新日期：2025.02.01 13:16:44"
$ws.Range("D8").Value = "*maa://21476 (73.08), *maa://39431 (53.85), *maa://37551 (57.14)"
$ws.Range("H8").Value = "*maa://24371 (54.17)"
$ws.Range("D9").Value = "maa://22765 (92.55), *maa://21915 (70.37)"
$ws.Range("L9").Value = "maa://22762 (92.13), maa://39552 (81.82)"
$ws.Range("X9").Value = "maa://26223 (97.74)"
$ws.Range("AB9").Value = "maa://28711 (86.73), ***maa://22740 (5.77), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (95.45)"
$ws.Range("AF9").Value = "maa://26206 (90.43), *maa://22865 (50.94)"
$ws.Range("D10").Value = "***maa://25695 (18.82), **maa://32237 (41.3), ***maa://34206 (20.0), ***maa://39951 (14.89), ***maa://39243 (28.57), *maa://45271 (54.55)"
$ws.Range("T10").Value = "maa://27395 (96.34), maa://22755 (87.83), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("D11").Value = "maa://36707 (99.46)"
$ws.Range("T11").Value = "maa://22747 (92.9), maa://22501 (97.59), *maa://45521 (75.0)"
$ws.Range("AB11").Value = "maa://29912 (98.48), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("X12").Value = "maa://22753 (91.33), *maa://21485 (76.26), maa://37962 (88.89)"
$ws.Range("AB12").Value = "maa://23669 (95.47), maa://36677 (92.73), maa://39872 (90.91)"
$ws.Range("AF12").Value = "*maa://28932 (78.57), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (91.9), maa://36673 (93.15), maa://25001 (85.51)"
$ws.Range("T14").Value = "maa://22521 (94.12), maa://42751 (100.0)"
$ws.Range("D15").Value = "*maa://22743 (77.56), maa://22734 (84.03), *maa://30808 (64.18), **maa://36048 (42.11), maa://45058 (100.0)"
$ws.Range("T15").Value = "maa://23892 (97.44)"
$ws.Range("AF15").Value = "maa://21364 (81.37), *maa://36666 (78.57), *maa://22766 (68.97)"
$ws.Range("H18").Value = "maa://24421 (89.24)"
$ws.Range("L18").Value = "maa://22466 (89.61), *maa://22732 (51.16)"
$ws.Range("X18").Value = "maa://21917 (96.81), maa://22741 (85.71)"
$ws.Range("D20").Value = "maa://21432 (90.0), maa://25198 (93.4), *maa://20795 (51.56), maa://36680 (93.75)"
$ws.Range("L22").Value = "maa://27127 (80.87), *maa://22751 (72.06)"
$ws.Range("D23").Value = "***maa://28036 (28.57), *maa://41753 (53.33)"
$ws.Range("L23").Value = "maa://39756 (95.3), maa://39875 (94.12)"
$ws.Range("X23").Value = "*maa://28503 (65.75)"
$ws.Range("D24").Value = "*maa://24368 (77.81)"
$ws.Range("X24").Value = "maa://29988 (84.74), maa://23504 (93.19), **maa://22892 (40.14), *maa://25141 (76.74), *maa://36663 (78.08), ***maa://22815 (23.08)"
$ws.Range("H25").Value = "*maa://29063 (74.05), *maa://25311 (73.53), ***maa://22725 (4.84), *maa://45047 (71.43)"
$ws.Range("AB25").Value = "maa://31215 (87.27), *maa://24516 (80.0), maa://26001 (87.5)"
$ws.Range("AB26").Value = "maa://42235 (94.62)"
$ws.Range("AF28").Value = "maa://36660 (92.56), *maa://36701 (65.52)"
$ws.Range("L29").Value = "maa://28432 (93.25), *maa://28440 (79.05), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (68.77), maa://42865 (82.35), ***maa://34960 (8.33)"
$ws.Range("AB30").Value = "maa://42979 (96.48), maa://45822 (100.0), maa://45045 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.57), maa://36258 (85.32), *maa://43904 (72.73)"
$ws.Range("H32").Value = "maa://21895 (97.49), maa://36667 (98.61), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("L32").Value = "maa://28065 (95.35)"
$ws.Range("T32").Value = "maa://42859 (96.19), maa://41108 (88.0), maa://41238 (97.0), maa://45523 (100.0)"
$ws.Range("P33").Value = "maa://21956 (80.69), *maa://22730 (79.31)"
$ws.Range("H39").Value = "maa://36670 (88.3), maa://25199 (84.82), maa://30434 (90.79), ***maa://25036 (16.0), *maa://45059 (75.0), *maa://44165 (66.67)"
$ws.Range("P39").Value = "maa://24709 (91.67)"
$ws.Range("T39").Value = "maa://45788 (83.33), maa://45790 (88.89)"
$ws.Range("P41").Value = "**maa://35616 (38.24), maa://43177 (89.47)"
$ws.Range("H43").Value = "maa://22525 (92.31), maa://21284 (85.11)"
$ws.Range("H44").Value = "maa://29768 (97.95), maa://27728 (96.04)"
$ws.Range("H45").Value = "maa://21229 (84.49), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (84.21)"
$ws.Range("H46").Value = "maa://35931 (92.43), maa://43901 (88.24)"
$ws.Range("H55").Value = "maa://32532 (92.2)"
$ws.Range("H59").Value = "maa://27746 (82.88), maa://31270 (95.2)"
$ws.Range("H60").Value = "*maa://40438 (68.42)"
$ws.Range("H64").Value = "maa://44405 (88.46)"
